$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 08:10:38"
$ws1.Range("A3").Value = "Total filas: 107"

$ws1.Cells.Item(65, 1).Value = "07:12:47"
$ws1.Cells.Item(65, 2).Value = "07:36"
$ws1.Cells.Item(65, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(65, 4).Value = 24
$ws1.Cells.Item(65, 5).Value = "LP1912"
$ws1.Cells.Item(66, 1).Value = "06:46:37"
$ws1.Cells.Item(66, 2).Value = "07:36"
$ws1.Cells.Item(66, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(66, 4).Value = 50
$ws1.Cells.Item(66, 5).Value = "LP1912"
$ws1.Cells.Item(75, 1).Value = "07:50:33"
$ws1.Cells.Item(75, 2).Value = "07:59"
$ws1.Cells.Item(75, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(75, 4).Value = 9
$ws1.Cells.Item(75, 5).Value = "LP1912"
$ws1.Cells.Item(76, 1).Value = "06:53:56"
$ws1.Cells.Item(76, 2).Value = "07:59"
$ws1.Cells.Item(76, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(76, 4).Value = 66
$ws1.Cells.Item(76, 5).Value = "LP1912"
$ws1.Cells.Item(77, 1).Value = "06:46:37"
$ws1.Cells.Item(77, 2).Value = "08:00"
$ws1.Cells.Item(77, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(77, 4).Value = 74
$ws1.Cells.Item(77, 5).Value = "LP1912"
$ws1.Cells.Item(78, 1).Value = "06:18:01"
$ws1.Cells.Item(78, 2).Value = "08:00"
$ws1.Cells.Item(78, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(78, 4).Value = 102
$ws1.Cells.Item(78, 5).Value = "LP1912"
$ws1.Cells.Item(83, 1).Value = "08:10:38"
$ws1.Cells.Item(83, 2).Value = "08:19"
$ws1.Cells.Item(83, 3).Value = "17_ROMERO"
$ws1.Cells.Item(83, 4).Value = 9
$ws1.Cells.Item(83, 5).Value = "LP1912"
$ws1.Cells.Item(84, 1).Value = "08:10:38"
$ws1.Cells.Item(84, 2).Value = "08:21"
$ws1.Cells.Item(84, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(84, 4).Value = 11
$ws1.Cells.Item(84, 5).Value = "LP1912"
$ws1.Cells.Item(85, 1).Value = "08:10:38"
$ws1.Cells.Item(85, 2).Value = "08:29"
$ws1.Cells.Item(85, 3).Value = "14_ABASTO"
$ws1.Cells.Item(85, 4).Value = 19
$ws1.Cells.Item(85, 5).Value = "LP1912"
$ws1.Cells.Item(86, 1).Value = "08:10:38"
$ws1.Cells.Item(86, 2).Value = "08:33"
$ws1.Cells.Item(86, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(86, 4).Value = 23
$ws1.Cells.Item(86, 5).Value = "LP1912"
$ws1.Cells.Item(87, 1).Value = "08:10:38"
$ws1.Cells.Item(87, 2).Value = "08:33"
$ws1.Cells.Item(87, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(87, 4).Value = 23
$ws1.Cells.Item(87, 5).Value = "LP1912"
$ws1.Cells.Item(88, 1).Value = "07:50:33"
$ws1.Cells.Item(88, 2).Value = "08:34"
$ws1.Cells.Item(88, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(88, 4).Value = 44
$ws1.Cells.Item(88, 5).Value = "LP1912"
$ws1.Cells.Item(89, 1).Value = "08:10:38"
$ws1.Cells.Item(89, 2).Value = "08:41"
$ws1.Cells.Item(89, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(89, 4).Value = 31
$ws1.Cells.Item(89, 5).Value = "LP1912"
$ws1.Cells.Item(90, 1).Value = "07:12:47"
$ws1.Cells.Item(90, 2).Value = "08:47"
$ws1.Cells.Item(90, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(90, 4).Value = 95
$ws1.Cells.Item(90, 5).Value = "LP1912"
$ws1.Cells.Item(91, 1).Value = "08:10:38"
$ws1.Cells.Item(91, 2).Value = "08:48"
$ws1.Cells.Item(91, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(91, 4).Value = 38
$ws1.Cells.Item(91, 5).Value = "LP1912"
$ws1.Cells.Item(92, 1).Value = "08:10:38"
$ws1.Cells.Item(92, 2).Value = "08:48"
$ws1.Cells.Item(92, 3).Value = "10_OLMOS"
$ws1.Cells.Item(92, 4).Value = 38
$ws1.Cells.Item(92, 5).Value = "LP1912"
$ws1.Cells.Item(93, 1).Value = "08:10:38"
$ws1.Cells.Item(93, 2).Value = "08:51"
$ws1.Cells.Item(93, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(93, 4).Value = 41
$ws1.Cells.Item(93, 5).Value = "LP1912"
$ws1.Cells.Item(94, 1).Value = "08:10:38"
$ws1.Cells.Item(94, 2).Value = "08:59"
$ws1.Cells.Item(94, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(94, 4).Value = 49
$ws1.Cells.Item(94, 5).Value = "LP1912"
$ws1.Cells.Item(95, 1).Value = "08:10:38"
$ws1.Cells.Item(95, 2).Value = "09:00"
$ws1.Cells.Item(95, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(95, 4).Value = 50
$ws1.Cells.Item(95, 5).Value = "LP1912"
$ws1.Cells.Item(96, 1).Value = "07:38:30"
$ws1.Cells.Item(96, 2).Value = "09:02"
$ws1.Cells.Item(96, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(96, 4).Value = 84
$ws1.Cells.Item(96, 5).Value = "LP1912"
$ws1.Cells.Item(97, 1).Value = "08:10:38"
$ws1.Cells.Item(97, 2).Value = "09:03"
$ws1.Cells.Item(97, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(97, 4).Value = 53
$ws1.Cells.Item(97, 5).Value = "LP1912"
$ws1.Cells.Item(98, 1).Value = "07:50:33"
$ws1.Cells.Item(98, 2).Value = "09:03"
$ws1.Cells.Item(98, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(98, 4).Value = 73
$ws1.Cells.Item(98, 5).Value = "LP1912"
$ws1.Cells.Item(99, 1).Value = "08:10:38"
$ws1.Cells.Item(99, 2).Value = "09:10"
$ws1.Cells.Item(99, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(99, 4).Value = 60
$ws1.Cells.Item(99, 5).Value = "LP1912"
$ws1.Cells.Item(100, 1).Value = "07:50:33"
$ws1.Cells.Item(100, 2).Value = "09:12"
$ws1.Cells.Item(100, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(100, 4).Value = 82
$ws1.Cells.Item(100, 5).Value = "LP1912"
$ws1.Cells.Item(101, 1).Value = "08:10:38"
$ws1.Cells.Item(101, 2).Value = "09:14"
$ws1.Cells.Item(101, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(101, 4).Value = 64
$ws1.Cells.Item(101, 5).Value = "LP1912"
$ws1.Cells.Item(102, 1).Value = "07:38:30"
$ws1.Cells.Item(102, 2).Value = "09:15"
$ws1.Cells.Item(102, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(102, 4).Value = 97
$ws1.Cells.Item(102, 5).Value = "LP1912"
$ws1.Cells.Item(103, 1).Value = "08:10:38"
$ws1.Cells.Item(103, 2).Value = "09:18"
$ws1.Cells.Item(103, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(103, 4).Value = 68
$ws1.Cells.Item(103, 5).Value = "LP1912"
$ws1.Cells.Item(104, 1).Value = "07:50:33"
$ws1.Cells.Item(104, 2).Value = "09:19"
$ws1.Cells.Item(104, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(104, 4).Value = 89
$ws1.Cells.Item(104, 5).Value = "LP1912"
$ws1.Cells.Item(105, 1).Value = "08:10:38"
$ws1.Cells.Item(105, 2).Value = "09:29"
$ws1.Cells.Item(105, 3).Value = "10_OLMOS"
$ws1.Cells.Item(105, 4).Value = 79
$ws1.Cells.Item(105, 5).Value = "LP1912"
$ws1.Cells.Item(106, 1).Value = "08:10:38"
$ws1.Cells.Item(106, 2).Value = "09:34"
$ws1.Cells.Item(106, 3).Value = "15_ABASTO"
$ws1.Cells.Item(106, 4).Value = 84
$ws1.Cells.Item(106, 5).Value = "LP1912"
$ws1.Cells.Item(107, 1).Value = "08:10:38"
$ws1.Cells.Item(107, 2).Value = "09:44"
$ws1.Cells.Item(107, 3).Value = "14_ABASTO"
$ws1.Cells.Item(107, 4).Value = 94
$ws1.Cells.Item(107, 5).Value = "LP1912"
$ws1.Cells.Item(108, 1).Value = "08:10:38"
$ws1.Cells.Item(108, 2).Value = "09:49"
$ws1.Cells.Item(108, 3).Value = "15_ABASTO"
$ws1.Cells.Item(108, 4).Value = 99
$ws1.Cells.Item(108, 5).Value = "LP1912"
$ws1.Cells.Item(109, 1).Value = "08:10:38"
$ws1.Cells.Item(109, 2).Value = "09:51"
$ws1.Cells.Item(109, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(109, 4).Value = 101
$ws1.Cells.Item(109, 5).Value = "LP1912"
$ws1.Cells.Item(110, 1).Value = "08:10:38"
$ws1.Cells.Item(110, 2).Value = "09:56"
$ws1.Cells.Item(110, 3).Value = "10_OLMOS"
$ws1.Cells.Item(110, 4).Value = 106
$ws1.Cells.Item(110, 5).Value = "LP1912"
$ws1.Cells.Item(111, 1).Value = "08:10:38"
$ws1.Cells.Item(111, 2).Value = "10:03"
$ws1.Cells.Item(111, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(111, 4).Value = 113
$ws1.Cells.Item(111, 5).Value = "LP1912"
$ws1.Cells.Item(112, 1).Value = "08:10:38"
$ws1.Cells.Item(112, 2).Value = "10:08"
$ws1.Cells.Item(112, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(112, 4).Value = 118
$ws1.Cells.Item(112, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:10:38"
$ws2.Range("A3").Value = "Total filas: 20"

$ws2.Cells.Item(18, 1).Value = "08:10:38"
$ws2.Cells.Item(18, 2).Value = "08:33"
$ws2.Cells.Item(18, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(18, 4).Value = 23
$ws2.Cells.Item(18, 5).Value = "LP1912"
$ws2.Cells.Item(21, 1).Value = "08:10:38"
$ws2.Cells.Item(21, 2).Value = "08:48"
$ws2.Cells.Item(21, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(21, 4).Value = 38
$ws2.Cells.Item(21, 5).Value = "LP1912"
$ws2.Cells.Item(22, 1).Value = "08:10:38"
$ws2.Cells.Item(22, 2).Value = "08:59"
$ws2.Cells.Item(22, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(22, 4).Value = 49
$ws2.Cells.Item(22, 5).Value = "LP1912"
$ws2.Cells.Item(23, 1).Value = "08:10:38"
$ws2.Cells.Item(23, 2).Value = "09:18"
$ws2.Cells.Item(23, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(23, 4).Value = 68
$ws2.Cells.Item(23, 5).Value = "LP1912"
$ws2.Cells.Item(24, 1).Value = "07:50:33"
$ws2.Cells.Item(24, 2).Value = "09:19"
$ws2.Cells.Item(24, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(24, 4).Value = 89
$ws2.Cells.Item(24, 5).Value = "LP1912"
$ws2.Cells.Item(25, 1).Value = "08:10:38"
$ws2.Cells.Item(25, 2).Value = "10:03"
$ws2.Cells.Item(25, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(25, 4).Value = 113
$ws2.Cells.Item(25, 5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:10:38"
$ws3.Range("A3").Value = "Total filas: 11"

$ws3.Cells.Item(10, 1).Value = "08:10:38"
$ws3.Cells.Item(10, 2).Value = "08:12"
$ws3.Cells.Item(10, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 2
$ws3.Cells.Item(10, 5).Value = "L6173"
$ws3.Cells.Item(11, 1).Value = "07:12:47"
$ws3.Cells.Item(11, 2).Value = "08:22"
$ws3.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(11, 4).Value = 70
$ws3.Cells.Item(11, 5).Value = "L6203"
$ws3.Cells.Item(12, 1).Value = "07:50:33"
$ws3.Cells.Item(12, 2).Value = "08:23"
$ws3.Cells.Item(12, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(12, 4).Value = 33
$ws3.Cells.Item(12, 5).Value = "L6203"
$ws3.Cells.Item(13, 1).Value = "08:10:38"
$ws3.Cells.Item(13, 2).Value = "08:24"
$ws3.Cells.Item(13, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(13, 4).Value = 14
$ws3.Cells.Item(13, 5).Value = "L6203"
$ws3.Cells.Item(14, 1).Value = "08:10:38"
$ws3.Cells.Item(14, 2).Value = "08:51"
$ws3.Cells.Item(14, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(14, 4).Value = 41
$ws3.Cells.Item(14, 5).Value = "L6173"
$ws3.Cells.Item(15, 1).Value = "07:50:33"
$ws3.Cells.Item(15, 2).Value = "08:52"
$ws3.Cells.Item(15, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(15, 4).Value = 62
$ws3.Cells.Item(15, 5).Value = "L6173"
$ws3.Cells.Item(16, 1).Value = "08:10:38"
$ws3.Cells.Item(16, 2).Value = "09:55"
$ws3.Cells.Item(16, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(16, 4).Value = 105
$ws3.Cells.Item(16, 5).Value = "L6203"
